{"js": "// Replace each \"old\" three-digit-by-one-digit multiplication expression\n// with its corresponding \"new\" expression, one-to-one, inside the document\n// body (the values live in a 5-column table of practice problems).\nconst replacements = [\n  [\"505\u00d73=1515\", \"101\u00d73=303\"],\n  [\"552\u00d78=4416\", \"138\u00d73=414\"],\n  [\"286\u00d72=572\", \"288\u00d77=2016\"],\n  [\"757\u00d79=6813\", \"170\u00d75=850\"],\n  [\"795\u00d73=2385\", \"321\u00d79=2889\"],\n  [\"827\u00d78=6616\", \"587\u00d77=4109\"],\n  [\"255\u00d72=510\", \"169\u00d72=338\"],\n  [\"399\u00d78=3192\", \"464\u00d78=3712\"],\n  [\"223\u00d76=1338\", \"307\u00d75=1535\"],\n  [\"621\u00d72=1242\", \"757\u00d77=5299\"],\n  [\"209\u00d74=836\", \"452\u00d72=904\"],\n  [\"527\u00d72=1054\", \"384\u00d72=768\"],\n  [\"192\u00d75=960\", \"131\u00d73=393\"],\n  [\"699\u00d78=5592\", \"460\u00d72=920\"],\n  [\"510\u00d78=4080\", \"526\u00d74=2104\"],\n  [\"764\u00d75=3820\", \"545\u00d74=2180\"],\n  [\"414\u00d73=1242\", \"305\u00d77=2135\"],\n  [\"528\u00d77=3696\", \"173\u00d72=346\"],\n  [\"156\u00d77=1092\", \"134\u00d72=268\"],\n  [\"623\u00d74=2492\", \"255\u00d79=2295\"],\n  [\"197\u00d78=1576\", \"697\u00d72=1394\"],\n  [\"367\u00d78=2936\", \"650\u00d78=5200\"],\n  [\"453\u00d79=4077\", \"984\u00d74=3936\"],\n  [\"120\u00d76=720\", \"342\u00d75=1710\"],\n  [\"257\u00d76=1542\", \"378\u00d74=1512\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-by-one-digit multiplication expression with its\n# corresponding new expression across the document (practice-problem table).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"505\u00d73=1515\", \"101\u00d73=303\"),\n    @(\"552\u00d78=4416\", \"138\u00d73=414\"),\n    @(\"286\u00d72=572\", \"288\u00d77=2016\"),\n    @(\"757\u00d79=6813\", \"170\u00d75=850\"),\n    @(\"795\u00d73=2385\", \"321\u00d79=2889\"),\n    @(\"827\u00d78=6616\", \"587\u00d77=4109\"),\n    @(\"255\u00d72=510\", \"169\u00d72=338\"),\n    @(\"399\u00d78=3192\", \"464\u00d78=3712\"),\n    @(\"223\u00d76=1338\", \"307\u00d75=1535\"),\n    @(\"621\u00d72=1242\", \"757\u00d77=5299\"),\n    @(\"209\u00d74=836\", \"452\u00d72=904\"),\n    @(\"527\u00d72=1054\", \"384\u00d72=768\"),\n    @(\"192\u00d75=960\", \"131\u00d73=393\"),\n    @(\"699\u00d78=5592\", \"460\u00d72=920\"),\n    @(\"510\u00d78=4080\", \"526\u00d74=2104\"),\n    @(\"764\u00d75=3820\", \"545\u00d74=2180\"),\n    @(\"414\u00d73=1242\", \"305\u00d77=2135\"),\n    @(\"528\u00d77=3696\", \"173\u00d72=346\"),\n    @(\"156\u00d77=1092\", \"134\u00d72=268\"),\n    @(\"623\u00d74=2492\", \"255\u00d79=2295\"),\n    @(\"197\u00d78=1576\", \"697\u00d72=1394\"),\n    @(\"367\u00d78=2936\", \"650\u00d78=5200\"),\n    @(\"453\u00d79=4077\", \"984\u00d74=3936\"),\n    @(\"120\u00d76=720\", \"342\u00d75=1710\"),\n    @(\"257\u00d76=1542\", \"378\u00d74=1512\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n\n"}
